$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# AnalogIO(3) sheet (sheet7.xml) - implement multi-sample analog IO
# -----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AnalogIO(3)")

# --- Command section (rows 1-7) ---
# Row 6: "Samples"/uint8 -> "MinSampleInterval"/uint16, with a new note.
$ws.Range("B6").Value = "MinSampleInterval"
$ws.Range("C6").Value = "uint16"
$ws.Range("F6").Value = "Minimum delay between samples, may not be slower than communication reate"

# --- Status section (rows 9-17) ---
# Row 14: "Value"/uint16/len2 -> "Samples"/uint8/len1
$ws.Range("B14").Value = "Samples"
$ws.Range("C14").Value = "uint8"
$ws.Range("D14").Value = 1

# Row 15: new sub-header row (merged B15:F15), styled like the other
# "AnalogIOConfig[AnalogConfigCount]" sub-header rows (bold, bottom border,
# left aligned - same look as B4/B12).
$ws.Range("B15:F15").Value = ""
$ws.Range("B4:F4").Copy()
$ws.Range("B15:F15").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B15").Value = "Samples[Samples]"
$ws.Range("B15").Characters(9, 7).Font.Italic = $true
$ws.Range("C15").Value = "AIOSample"
$ws.Range("B15:F15").Merge()

# Row 16: new "Delay" field (uint32, offset resets to 0 inside the struct)
$ws.Range("B16").Value = "Delay"
$ws.Range("C16").Value = "uint32"
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = "Delay since previous sample"

# Row 17: "Value" field (uint16), moved down from the old row 14/15
$ws.Range("B17").Value = "Value"
$ws.Range("C17").Value = "uint16"
$ws.Range("D17").Value = 2
$ws.Range("E17").Formula = "=E16+D16"

# Column F needs to be widened to fit the long new notes text.
$ws.Columns.Item(6).ColumnWidth = 74.140625

# Make AnalogIO(3) the active/selected sheet with a particular selection,
# matching the updated workbook view.
$ws.Activate()
$ws.Range("B19").Select()

# -----------------------------------------------------------------------
# Sections sheet (sheet4.xml) - no longer the active tab; selection moves.
# -----------------------------------------------------------------------
$wsSections = $wb.Worksheets.Item("Sections")
$wsSections.Activate()
$wsSections.Range("B8").Select()

# Re-activate AnalogIO(3) last so it ends up as the selected/active tab.
$ws.Activate()
$ws.Range("B19").Select()
